{"js": "// Locate the \"Research Fields\" paragraph text and apply the two edits:\n//  1. Insert \" Racial Disparity,\" right after \"Health Economics,\"\n//  2. Remove the trailing \", Quantile.\" after \"...Nonparametric\"\n//\n// Both target substrings are unique within the document, so a body-wide\n// search is sufficient to anchor each edit.\n\nconst insertAfter = context.document.body.search(\"Health Economics,\", { matchCase: true });\ninsertAfter.load(\"items\");\nawait context.sync();\n\nif (insertAfter.items.length > 0) {\n  insertAfter.items[0].insertText(\" Racial Disparity,\", \"After\");\n  await context.sync();\n}\n\nconst toRemove = context.document.body.search(\", Quantile.\", { matchCase: true });\ntoRemove.load(\"items\");\nawait context.sync();\n\nif (toRemove.items.length > 0) {\n  toRemove.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop script: add \"Racial Disparity,\" to the Research Fields\n# list, and drop the trailing \", Quantile.\" after \"Nonparametric\".\n$d = $word.ActiveDocument\n\n# 1) \"Applied Econometrics, Health Economics, Casusal ...\" ->\n#    \"Applied Econometrics, Health Economics, Racial Disparity, Casusal ...\"\n$find1 = $d.Content.Find\n$find1.Text = \"Health Economics,\"\n$find1.Replacement.Text = \"Health Economics, Racial Disparity,\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) \"...Nonparametric, Quantile.\" -> \"...Nonparametric\"\n$find2 = $d.Content.Find\n$find2.Text = \", Quantile.\"\n$find2.Replacement.Text = \"\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
